$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.180.72'
$ws.Range("D3").Value = '1.786.78'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").Value = '''226.22'
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("D6").Value = '''0.547'
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").Value = '''31.96'
$ws.Range("E8").Value = '  -2.37%  '
$ws.Range("D9").Value = '''0.292'
$ws.Range("E9").Value = '  +2.25%  '
$ws.Range("D10").Value = '''0.0690'
$ws.Range("E10").Value = '  -3.03%  '
$ws.Range("D11").Value = '''0.0944'
$ws.Range("E11").Value = '  +1.02%  '
$ws.Range("D12").Value = '2.045.50'
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '''11.23'
$ws.Range("E13").Value = '  +0.53%  '
$ws.Range("D14").Value = '1.795.46'
$ws.Range("E14").Value = '  +0.76%  '
$ws.Range("D15").Value = '34.121.45'
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("D16").Value = '''0.620'
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("E17").Value = '  +1.87%  '
$ws.Range("D18").Value = '''68.00'
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("D19").Value = '''246.47'
$ws.Range("E19").Value = '  +0.91%  '
$ws.Range("D20").Value = '0.0₃0780'
$ws.Range("E20").Value = '  -0.36%  '
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("D22").Value = '''10.83'
$ws.Range("E22").Value = '  +1.18%  '
$ws.Range("E23").Value = '  +0.68%  '
$ws.Range("E24").Value = '  -1.32%  '
$ws.Range("D25").Value = '''161.46'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("E26").Value = '  +1.18%  '
$ws.Range("D27").Value = '''16.34'
$ws.Range("E27").Value = '  +0.18%  '
$ws.Range("E28").Value = '  +1.49%  '
$ws.Range("E29").Value = '  +0.33%  '
$ws.Range("D30").Value = '''1.23'
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("E31").Value = '  +1.67%  '
$ws.Range("E32").Value = '  +0.66%  '
$ws.Range("E33").Value = '  +3.30%  '
$ws.Range("D34").Value = '''1.81'
$ws.Range("E34").Value = '  +0.85%  '
$ws.Range("D35").Value = '1.449.79'
$ws.Range("E35").Value = '  +4.18%  '
$ws.Range("D36").Value = '''0.646'
$ws.Range("E36").Value = '  -0.66%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '''0.0193'
$ws.Range("E37").Value = '  +3.39%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '''2.41'
$ws.Range("E38").Value = '  +9.45%  '
$ws.Range("E39").Value = '  -0.61%  '
$ws.Range("D40").Value = '''80.03'
$ws.Range("E40").Value = '  +2.68%  '
$ws.Range("D41").Value = '''2.37'
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("D42").Value = '''0.920'
$ws.Range("E42").Value = '  +0.79%  '
$ws.Range("D43").Value = '''2.68'
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").Value = '''13.40'
$ws.Range("E44").Value = '  +2.03%  '
$ws.Range("B45").Value = 'Kaspa'
$ws.Range("C45").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D45").Value = '''0.0510'
$ws.Range("E45").Value = '  +2.74%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '''6.07'
$ws.Range("E46").Value = '  +5.09%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₆0138'
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("D49").Value = '''107.77'
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("D50").Value = '1.946.84'
$ws.Range("E50").Value = '  +0.61%  '
$ws.Range("E51").Value = '  +0.18%  '
